# "Rounding measurement + adding sensors"
#
# The sensor-ID column (A) holds long numeric-looking strings (e.g.
# "28072261300627"). Typed in as plain values they would be interpreted
# as numbers and silently rounded/reformatted by Excel, so the column is
# switched to a Text number format before the new IDs are entered, and
# four new temperature sensors (with both circuit-labelled and generic
# "Aanvoer"/"Afvoer" names) are appended to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing ID-column formatting (Courier New font, from the
# table's "ID" column style) down onto the rows that are about to receive
# new sensor IDs, so the whole column keeps one consistent style.
$ws.Range("A6").Copy()
$ws.Range("A7:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Store the ID column as Text so the long numeric IDs aren't rounded.
$ws.Range("A2:A11").NumberFormat = "@"

# Name the four sensors that already had an ID but no label yet.
$ws.Range("B2").Value = "Kring 1 aanvoer"
$ws.Range("B3").Value = "Kring 1 afvoer"
$ws.Range("B4").Value = "Kring 2 aanvoer"
$ws.Range("B5").Value = "Kring 2 afvoer"

# The sensor that used to be the last row moves down to row 10, relabelled.
$ws.Range("A10").Value = "282bfe571f64ff"
$ws.Range("B10").Value = "Aanvoer"

# New sensor IDs for the new circuits.
$ws.Range("A6").Value = "28072261300627"
$ws.Range("A7").Value = "280722613294cc"
$ws.Range("A8").Value = "280722614c7990"
$ws.Range("A9").Value = "280922545d1f8a"
$ws.Range("A11").Value = "28092254776424"

# Their labels.
$ws.Range("B6").Value = "Kring 3 aanvoer"
$ws.Range("B7").Value = "Kring 3 afvoer"
$ws.Range("B8").Value = "Kring 4 aanvoer"
$ws.Range("B9").Value = "Kring 4 afvoer"

$ws.Range("B11").Value = "Afvoer"

# Grow the table (and its autofilter) to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B11"))

$ws.Range("A7").Select() | Out-Null
